$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.97'
$ws.Range("G2").Value = '''15'

$ws.Range("D3").Value = '''26.26'
$ws.Range("G3").Value = '''15'

$ws.Range("D4").Value = '''5.074'
$ws.Range("G4").Value = '''15'

$ws.Range("D5").Value = '''0.05615'
$ws.Range("G5").Value = '''15'

$ws.Range("D6").Value = '''6.513'
$ws.Range("G6").Value = '''15'

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''3.058'
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("G7").Value = '''15'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.8132'
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("G8").Value = '''15'

$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = '''0.8498'
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("G9").Value = '''15'

$ws.Range("D10").Value = '''0.02823'
$ws.Range("G10").Value = '''15'

$ws.Range("D11").Value = '''0.09398'
$ws.Range("G11").Value = '''15'

$ws.Range("D12").Value = '''0.001520'
$ws.Range("G12").Value = '''15'

$ws.Range("D13").Value = '''0.0005971'
$ws.Range("G13").Value = '''15'

$ws.Range("D14").Value = '''0.006150'
$ws.Range("G14").Value = '''15'

$ws.Range("D15").Value = '''3.578'
$ws.Range("G15").Value = '''15'

$ws.Range("B16").Value = 'BTSEToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D16").Value = '''2.118'
$ws.Range("E16").Value = '15BTSETokenBTSE'
$ws.Range("G16").Value = '''15'

$ws.Range("B17").Value = 'BitpandaEcosystemToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D17").Value = '''0.3181'
$ws.Range("E17").Value = '16BitpandaEcosystemTokenBEST'
$ws.Range("G17").Value = '''15'

$ws.Range("B18").Value = 'WazirX'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D18").Value = '''0.1340'
$ws.Range("E18").Value = '17WazirXWRX'
$ws.Range("G18").Value = '''15'

$ws.Range("B19").Value = 'MandalaExchangeToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D19").Value = '''0.06966'
$ws.Range("E19").Value = '18MandalaExchangeTokenMDX'
$ws.Range("G19").Value = '''15'

$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '''0.03166'
$ws.Range("E20").Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G20").Value = '''15'

$ws.Range("D21").Value = '''0.1321'
$ws.Range("G21").Value = '''15'

$ws.Range("D22").Value = '''3.739'
$ws.Range("G22").Value = '''15'

$ws.Range("D23").Value = '''0.04672'
$ws.Range("G23").Value = '''15'

$ws.Range("D24").Value = '''0.1376'
$ws.Range("G24").Value = '''15'

$ws.Range("D25").Value = '''0.001250'
$ws.Range("G25").Value = '''15'

$ws.Range("D26").Value = '''0.004620'
$ws.Range("G26").Value = '''15'

$ws.Range("D27").Value = '''0.00009599'
$ws.Range("G27").Value = '''15'

$ws.Range("D28").Value = '''0.0001938'
$ws.Range("G28").Value = '''15'

$ws.Range("G29").Value = '''15'

$ws.Range("G30").Value = '''15'

$ws.Range("G31").Value = '''15'

$ws.Range("G32").Value = '''15'

$ws.Range("G33").Value = '''15'

$ws.Range("G34").Value = '''15'

$ws.Range("G35").Value = '''15'

$ws.Range("G36").Value = '''15'

$ws.Range("G37").Value = '''15'

$ws.Range("G38").Value = '''15'

$ws.Range("G39").Value = '''15'

$ws.Range("D40").Value = '''0.03665'
$ws.Range("G40").Value = '''15'

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1364'
$ws.Range("E41").Value = '40BKEXTokenBKKBestin24h'
$ws.Range("G41").Value = '''15'

$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = '''0.006089'
$ws.Range("E42").Value = '41KickTokenKICK'
$ws.Range("G42").Value = '''15'

$ws.Range("D43").Value = '''0.002660'
$ws.Range("G43").Value = '''15'

$ws.Range("D44").Value = '''0.008698'
$ws.Range("G44").Value = '''15'

$ws.Range("D45").Value = '''0.00005292'
$ws.Range("G45").Value = '''15'

$ws.Range("G46").Value = '''15'

$ws.Range("G47").Value = '''15'

$ws.Range("D48").Value = '''0.002070'
$ws.Range("G48").Value = '''15'

$ws.Range("D49").Value = '''0.00002100'
$ws.Range("G49").Value = '''15'

$ws.Range("D50").Value = '''0.0002000'
$ws.Range("G50").Value = '''15'

$ws.Range("G51").Value = '''15'

